$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.760.91"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "3.274.13"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'584.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'182.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "'0.423"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "3.854.85"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'28.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "68.806.42"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "3.255.48"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'5.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'13.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'394.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.92%  "
$ws.Range("D21").Value = "'7.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "'71.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").Value = "'0.0000120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +3.77%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "'1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "'5.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "'22.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").Value = "'7.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D35").Value = "'164.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "'0.830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").Value = "'4.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -3.68%  "
$ws.Range("D42").Value = "'2.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "'41.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "'0.0687"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").Value = "'344.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.67%  "
$ws.Range("D46").Value = "2.615.33"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").Value = "'24.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'0.0281"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "'31.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'6.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +0.16%  "
